# Apply the diff: swap the match data stored in rows 26 & 28 (cols F:V),
# and append two new match rows (59 & 60) with updated "used range".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row 26 <-> Row 28 swap (only columns F..V - A..E index/meta untouched)
# ---------------------------------------------------------------------

# Values currently in row 26 (Kolding IF vs Sonderjyske match)
$row26 = @("Kolding IF", 1, "Sonderjyske", 3, 2.08, "14/08/2023 04:12", 2.99, "18/08/2023 18:30", 3.6, "14/08/2023 04:12", 3.45, "18/08/2023 18:31", 3.38, "14/08/2023 04:12", 2.38, "18/08/2023 18:30", "https://www.betexplorer.com/football/denmark/1st-division/kolding-if-sonderjyske/UiWf7qoG/")

# Values currently in row 28 (Horsens vs Helsingor match)
$row28 = @("Horsens", 3, "Helsingor", 1, 2.01, "13/08/2023 22:12", 2.31, "18/08/2023 18:51", 3.79, "13/08/2023 22:12", 3.76, "18/08/2023 18:51", 3.2, "13/08/2023 22:12", 2.89, "18/08/2023 18:51", "https://www.betexplorer.com/football/denmark/1st-division/horsens-helsingor/nTtj8PV9/")

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "26").Value = $row28[$i]
    $ws.Range($cols[$i] + "28").Value = $row26[$i]
}

# ---------------------------------------------------------------------
# 2) Append new rows 59 and 60
# ---------------------------------------------------------------------

$ws.Range("A59").Value = 58
$ws.Range("B59").Value = "denmark"
$ws.Range("C59").Value = "1st-division"
$ws.Range("D59").Value = "2023-2024"
$ws.Range("E59").Value = 45193.54166666666
$ws.Range("F59").Value = "Naestved"
$ws.Range("G59").Value = 2
$ws.Range("H59").Value = "Hillerod"
$ws.Range("I59").Value = 2
$ws.Range("J59").Value = 2.24
$ws.Range("K59").Value = "20/09/2023 00:42"
$ws.Range("L59").Value = 2.66
$ws.Range("M59").Value = "24/09/2023 07:41"
$ws.Range("N59").Value = 3.76
$ws.Range("O59").Value = "20/09/2023 00:42"
$ws.Range("P59").Value = 3.7
$ws.Range("Q59").Value = "24/09/2023 12:50"
$ws.Range("R59").Value = 3.01
$ws.Range("S59").Value = "20/09/2023 00:42"
$ws.Range("T59").Value = 2.51
$ws.Range("U59").Value = "24/09/2023 12:50"
$ws.Range("V59").Value = "https://www.betexplorer.com/football/denmark/1st-division/naestved-if-hillerod/rJzKzmSk/"

$ws.Range("A60").Value = 59
$ws.Range("B60").Value = "denmark"
$ws.Range("C60").Value = "1st-division"
$ws.Range("D60").Value = "2023-2024"
$ws.Range("E60").Value = 45193.58333333334
$ws.Range("F60").Value = "Helsingor"
$ws.Range("G60").Value = 3
$ws.Range("H60").Value = "Sonderjyske"
$ws.Range("I60").Value = 4
$ws.Range("J60").Value = 3.56
$ws.Range("K60").Value = "17/09/2023 13:12"
$ws.Range("L60").Value = 4.96
$ws.Range("M60").Value = "24/09/2023 13:51"
$ws.Range("N60").Value = 3.89
$ws.Range("O60").Value = "17/09/2023 13:12"
$ws.Range("P60").Value = 4.58
$ws.Range("Q60").Value = "24/09/2023 13:51"
$ws.Range("R60").Value = 1.93
$ws.Range("S60").Value = "17/09/2023 13:12"
$ws.Range("T60").Value = 1.59
$ws.Range("U60").Value = "24/09/2023 13:51"
$ws.Range("V60").Value = "https://www.betexplorer.com/football/denmark/1st-division/helsingor-sonderjyske/jLeurVZF/"

# ---------------------------------------------------------------------
# 3) Match the formatting used by the other data rows:
#    column A -> bordered/bold "index" style, column E -> date-time format
# ---------------------------------------------------------------------

$ws.Range("A58").Copy()
$ws.Range("A59").PasteSpecial(-4122)
$ws.Range("A60").PasteSpecial(-4122)

$ws.Range("E58").Copy()
$ws.Range("E59").PasteSpecial(-4122)
$ws.Range("E60").PasteSpecial(-4122)
